$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the test-case data: update the two claim numbers (NroSiniestro)
# and the tester user name that were wrong.
# F2/F3 already carry the "quote-prefixed text" style (leading apostrophe),
# so prefix the new values the same way to keep them as plain text without
# altering the cell's NumberFormat/style.
$ws.Range("F3").Value = "'0420172010222  "
$ws.Range("F2").Value = "'0420194406906"
$ws.Range("D2").Value = "tcorvetto"

# Move the active selection to D3 (cosmetic - matches saved cursor position)
$ws.Range("D3").Select()
